$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.965722762100373
$ws.Range("C2").Value = 5.141302184793021
$ws.Range("D2").Value = 5.033115134359021
$ws.Range("E2").Value = 12.80736667787291
$ws.Range("F2").Value = 24.52540140325567
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 20.88698452704859
$ws.Range("K2").Value = 8.019165173221673
$ws.Range("M2").Value = 13.34531255437553
$ws.Range("O2").Value = 21.98715207275361
$ws.Range("B3").Value = 7.654023975218778
$ws.Range("C3").Value = 4.955890905370541
$ws.Range("D3").Value = 4.977371347440244
$ws.Range("E3").Value = 12.60303087404094
$ws.Range("F3").Value = 24.55218966510247
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 20.97441576261757
$ws.Range("K3").Value = 7.732999789525001
$ws.Range("M3").Value = 13.1669280627894
$ws.Range("O3").Value = 22.05587897094495
$ws.Range("B4").Value = 7.456945490711633
$ws.Range("C4").Value = 4.837442103993507
$ws.Range("D4").Value = 4.94231755823872
$ws.Range("E4").Value = 12.48047349263396
$ws.Range("F4").Value = 24.57542221950061
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 21.03228952619626
$ws.Range("K4").Value = 7.550151859790227
$ws.Range("M4").Value = 13.05933914864445
$ws.Range("O4").Value = 22.10305660923381
$ws.Range("B5").Value = 7.375323452522463
$ws.Range("C5").Value = 4.788058318119004
$ws.Range("D5").Value = 4.927831937099579
$ws.Range("E5").Value = 12.43133169749124
$ws.Range("F5").Value = 24.58659232514336
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 21.05692594679654
$ws.Range("K5").Value = 7.476200280887971
$ws.Range("M5").Value = 13.01603650315322
$ws.Range("O5").Value = 22.12353069553694
$ws.Range("B6").Value = 7.361694797733367
$ws.Range("C6").Value = 4.779792200726477
$ws.Range("D6").Value = 4.925414683452067
$ws.Range("E6").Value = 12.4232221718813
$ws.Range("F6").Value = 24.58854983368147
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 21.06108032764654
$ws.Range("K6").Value = 7.467170679401525
$ws.Range("M6").Value = 13.00888031663797
$ws.Range("O6").Value = 22.12700573098022
$ws.Range("B7").Value = 7.455849845605031
$ws.Range("C7").Value = 4.836780550537441
$ws.Range("D7").Value = 4.942123004239892
$ws.Range("E7").Value = 12.47980741243556
$ws.Range("F7").Value = 24.57556597457791
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 21.03261752230959
$ws.Range("K7").Value = 7.549130545273555
$ws.Range("M7").Value = 13.05875289551945
$ws.Range("O7").Value = 22.10332767787502
$ws.Range("B8").Value = 7.85949920118712
$ws.Range("C8").Value = 5.078352582384415
$ws.Range("D8").Value = 5.014070449092886
$ws.Range("E8").Value = 12.73635148813495
$ws.Range("F8").Value = 24.53322811958222
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 20.91626022354453
$ws.Range("K8").Value = 7.922014561031916
$ws.Range("M8").Value = 13.28343496710814
$ws.Range("O8").Value = 22.00981381961502
$ws.Range("B9").Value = 8.601234424360607
$ws.Range("C9").Value = 5.513886891620851
$ws.Range("D9").Value = 5.148299990107659
$ws.Range("E9").Value = 13.25932287702096
$ws.Range("F9").Value = 24.50414849026399
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 20.72139623791922
$ws.Range("K9").Value = 8.594070039828432
$ws.Range("M9").Value = 13.73706833516171
$ws.Range("O9").Value = 21.86607813232402
$ws.Range("B10").Value = 9.110477561380435
$ws.Range("C10").Value = 5.808712712280439
$ws.Range("D10").Value = 5.24235191485486
$ws.Range("E10").Value = 13.65140948460387
$ws.Range("F10").Value = 24.51576988589902
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 20.59862204425933
$ws.Range("K10").Value = 9.048888468086567
$ws.Range("M10").Value = 14.07512027044107
$ws.Range("O10").Value = 21.78482519913877
$ws.Range("B11").Value = 9.33349042001338
$ws.Range("C11").Value = 5.937057401422086
$ws.Range("D11").Value = 5.284069260736405
$ws.Range("E11").Value = 13.83059310120407
$ws.Range("F11").Value = 24.52822101173205
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 20.54721539386346
$ws.Range("K11").Value = 9.246858962109867
$ws.Range("M11").Value = 14.22927561087032
$ws.Range("O11").Value = 21.7531840052971
$ws.Range("B12").Value = 9.416634733166649
$ws.Range("C12").Value = 5.984806324020222
$ws.Range("D12").Value = 5.299706755941391
$ws.Range("E12").Value = 13.8984909589647
$ws.Range("F12").Value = 24.53396456655469
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 20.52838970739781
$ws.Range("K12").Value = 9.320508049894579
$ws.Range("M12").Value = 14.28764930269593
$ws.Range("O12").Value = 21.74197012221238
$ws.Range("B13").Value = 9.3987871006254
$ws.Range("C13").Value = 5.974560972415772
$ws.Range("D13").Value = 5.296346156374552
$ws.Range("E13").Value = 13.88386708427758
$ws.Range("F13").Value = 24.53268188013468
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 20.53241562197859
$ws.Range("K13").Value = 9.304705513978819
$ws.Range("M13").Value = 14.27507842103506
$ws.Range("O13").Value = 21.74435103576289
$ws.Range("B14").Value = 9.340357251516494
$ws.Range("C14").Value = 5.941002934310095
$ws.Range("D14").Value = 5.285359008945608
$ws.Range("E14").Value = 13.8361786099712
$ws.Range("F14").Value = 24.52867293627394
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 20.54565373807727
$ws.Range("K14").Value = 9.252944710849183
$ws.Range("M14").Value = 14.2340783936486
$ws.Range("O14").Value = 21.75224602635088
$ws.Range("B15").Value = 9.304395513631739
$ws.Range("C15").Value = 5.920336048848339
$ws.Range("D15").Value = 5.278608044382687
$ws.Range("E15").Value = 13.80697175704942
$ws.Range("F15").Value = 24.52635122747274
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 20.5538459855503
$ws.Range("K15").Value = 9.221067145753718
$ws.Range("M15").Value = 14.20896286649525
$ws.Range("O15").Value = 21.75718202471987
$ws.Range("B16").Value = 9.095723630263338
$ws.Range("C16").Value = 5.800206940111624
$ws.Range("D16").Value = 5.239603551681982
$ws.Range("E16").Value = 13.63970963483619
$ws.Range("F16").Value = 24.51510022490796
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 20.60207117220594
$ws.Range("K16").Value = 9.035767936540422
$ws.Range("M16").Value = 14.06504861764753
$ws.Range("O16").Value = 21.787000342659
$ws.Range("B17").Value = 8.96544862930509
$ws.Range("C17").Value = 5.725016138413457
$ws.Range("D17").Value = 5.215397506003725
$ws.Range("E17").Value = 13.53725736105801
$ws.Range("F17").Value = 24.51003224199992
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 20.63279528418384
$ws.Range("K17").Value = 8.919780200733564
$ws.Range("M17").Value = 13.97681784869511
$ws.Range("O17").Value = 21.8066578298347
$ws.Range("B18").Value = 8.889707392373952
$ws.Range("C18").Value = 5.68122610695287
$ws.Range("D18").Value = 5.201374738014434
$ws.Range("E18").Value = 13.47841109920351
$ws.Range("F18").Value = 24.50779179285754
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 20.65088515367997
$ws.Range("K18").Value = 8.852228425000545
$ws.Range("M18").Value = 13.9261082682991
$ws.Range("O18").Value = 21.81846499689932
$ws.Range("B19").Value = 8.863925565900804
$ws.Range("C19").Value = 5.666307150132978
$ws.Range("D19").Value = 5.196609880101023
$ws.Range("E19").Value = 13.45850304328195
$ws.Range("F19").Value = 24.50714909866644
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 20.65708183776954
$ws.Range("K19").Value = 8.829213626259497
$ws.Range("M19").Value = 13.90894718262748
$ws.Range("O19").Value = 21.82254860889826
$ws.Range("B20").Value = 8.979400992260974
$ws.Range("C20").Value = 5.733076621983919
$ws.Range("D20").Value = 5.217984687531545
$ws.Range("E20").Value = 13.54815567507469
$ws.Range("F20").Value = 24.51050193583576
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 20.62948135721955
$ws.Range("K20").Value = 8.932214357619424
$ws.Range("M20").Value = 13.98620659300732
$ws.Range("O20").Value = 21.80451341542687
$ws.Range("B21").Value = 9.357555410732175
$ws.Range("C21").Value = 5.950883040073371
$ws.Range("D21").Value = 5.288590592153175
$ws.Range("E21").Value = 13.8501852241268
$ws.Range("F21").Value = 24.52982256314471
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 20.54174797383128
$ws.Range("K21").Value = 9.26818413291684
$ws.Range("M21").Value = 14.24612157299832
$ws.Range("O21").Value = 21.74990621318276
$ws.Range("B22").Value = 9.597066772581899
$ws.Range("C22").Value = 6.088255351276716
$ws.Range("D22").Value = 5.333800434807976
$ws.Range("E22").Value = 14.04780851368942
$ws.Range("F22").Value = 24.54844358428128
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 20.48814564366165
$ws.Range("K22").Value = 9.480064684582491
$ws.Range("M22").Value = 14.41595712096998
$ws.Range("O22").Value = 21.71869459519734
$ws.Range("B23").Value = 9.469952109593855
$ws.Range("C23").Value = 6.015399099753579
$ws.Range("D23").Value = 5.309758740456101
$ws.Range("E23").Value = 13.94233582018649
$ws.Range("F23").Value = 24.53795757022237
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 20.51641170678642
$ws.Range("K23").Value = 9.367694209941261
$ws.Range("M23").Value = 14.32533371420367
$ws.Range("O23").Value = 21.7349422819033
$ws.Range("B24").Value = 8.973095754246845
$ws.Range("C24").Value = 5.729434225673981
$ws.Range("D24").Value = 5.216815353643788
$ws.Range("E24").Value = 13.54322837066218
$ws.Range("F24").Value = 24.51028749014061
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 20.63097825622491
$ws.Range("K24").Value = 8.926595579064315
$ws.Range("M24").Value = 13.98196188929646
$ws.Range("O24").Value = 21.8054813297787
$ws.Range("B25").Value = 8.406505697084169
$ws.Range("C25").Value = 5.400362215860284
$ws.Range("D25").Value = 5.112762138202227
$ws.Range("E25").Value = 13.11615760077366
$ws.Range("F25").Value = 24.50622602851738
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 20.77053715954398
$ws.Range("K25").Value = 8.418916833089551
$ws.Range("M25").Value = 13.61328190770182
$ws.Range("O25").Value = 21.90069884956073
